$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 186, pushing the existing rows 186-208 down to 187-209.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new weekly data point
# (same market/category metadata, new date + price observation).
$ws.Cells.Item(186, 1).Value = 10
$ws.Cells.Item(186, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(186, 3).Value = "La Araucanía"
$ws.Cells.Item(186, 4).Value = 44491
$ws.Cells.Item(186, 5).Value = 9
$ws.Cells.Item(186, 6).Value = 100112009
$ws.Cells.Item(186, 7).Value = "Acelga"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 65
$ws.Cells.Item(186, 11).Value = 7000
$ws.Cells.Item(186, 12).Value = 7000
$ws.Cells.Item(186, 13).Value = 7000
$ws.Cells.Item(186, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(186, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(186, 16).Value = 583
$ws.Cells.Item(186, 17).Value = 12
$ws.Cells.Item(186, 18).Value = "Hortaliza"
